$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LogInPage")
$ws2 = $wb.Worksheets.Item("ErrorMessages")

# --- ErrorMessages sheet: row 3 (BlankUserName) ---
$ws2.Range("A3").Value = "BlankUserName"
$ws2.Range("B3").Value = "Please enter Email / Mobile Number."
$ws2.Range("B2").Copy()
$ws2.Range("B3").PasteSpecial(-4122)

# --- ErrorMessages sheet: header B1 changes to Message ---
$ws2.Range("B1").Value = "Message"

# --- LogInPage sheet: C3 / C2 (InvalidValue / InvalidVallue) ---
$ws1.Range("B2").Copy()
$ws1.Range("C3").PasteSpecial(-4122)
$ws1.Range("C3").Value = "InvalidValue"

$ws1.Range("B2").Copy()
$ws1.Range("C2").PasteSpecial(-4122)
$ws1.Range("C2").Value = "InvalidVallue"

# --- ErrorMessages sheet: row 4 (IncorrectUserNameAndPassword) ---
$ws2.Range("A4").Value = "IncorrectUserNameAndPassword"
$ws2.Range("B4").Value = "Please check the username or password you have entered and also ensure you have selected the right profile – Buyer or Supplier"

# --- LogInPage sheet: C1 (InvalidValues header) ---
$ws1.Range("A1").Copy()
$ws1.Range("C1").PasteSpecial(-4122)
$ws1.Range("C1").Value = "InvalidValues"

# --- ErrorMessages sheet: row 5 (ForgotPasswordInvalidEmail) ---
$ws2.Range("B5").Value = "You have entered wrong Username / Mobile Number.Please try again or contact us at help@renepay.com."
$ws2.Range("A5").Value = "ForgotPasswordInvalidEmail"

# --- ErrorMessages sheet: row 6 (ForgotPasswordValidEmail) ---
$ws2.Range("B6").Value = "Password reset link has been sent on your registered email id.Please check your registered email account."
$ws2.Range("A6").Value = "ForgotPasswordValidEmail"

# --- ErrorMessages sheet: A2 LoginPage -> InvalidCredentials (reuses existing string) ---
$ws2.Range("A2").Value = "InvalidCredentials"

# --- LogInPage sheet: B3 picks up the centered style (same as B2/C2/C3) ---
$ws1.Range("B2").Copy()
$ws1.Range("B3").PasteSpecial(-4122)
$ws1.Range("B3").Value = "qwerty11"

# --- Column widths (values chosen so the engine's width quantization lands
#     as close as possible to the canonical widths 12.21875 / 28.44140625 / 107.21875) ---
$ws1.Columns.Item(3).ColumnWidth = 11.3
$ws2.Columns.Item(1).ColumnWidth = 27.65
$ws2.Columns.Item(2).ColumnWidth = 106.3

# --- Selections / active tab ---
$ws1.Range("A1:C3").Select()
$ws2.Range("B6").Select()
$ws2.Activate()
